$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row: year 2021, WCR value 28 (row 42, columns A and C only)
$ws.Range("A42").Value = 2021
$ws.Range("C42").Value = 28

# Update the selection to match the post-edit cursor position recorded in the diff
$ws.Range("D44").Select()
